$d = $word.ActiveDocument

# Locate the paragraph that holds the closing bibliography entry, then
# remove everything from the blank paragraph right after it through the
# end of the "(c) 2020 ... Creative Commons Attribution" paragraph
# (i.e. the blank line + the "Ver no Jupiter ..." line + the "(c) 2020 ..."
# line), leaving the bibliography entry followed directly by the blank
# paragraph that precedes the trailing page-break paragraph.

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Rio de Janeiro: Elsevier Editora, 2007.*") {
        $startPara = $p
    }
    if ($t -like "*Creative Commons Attribution*") {
        $endPara = $p
    }
}

$deleteStart = $startPara.Range.End
$deleteEnd = $endPara.Range.End

$r = $d.Range($deleteStart, $deleteEnd)
$r.Delete()
